$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6534544229507446
$ws.Range("B1").Value = 0.6806876659393311
$ws.Range("C1").Value = 2.116748332977295
$ws.Range("D1").Value = 3.648826599121094
$ws.Range("E1").Value = 1.256515264511108
